$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 87.25
$ws.Range("I2").Value = 87.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 87.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 25.75
$ws.Range("N2").ClearContents()

$ws.Range("H8").Value = 170.66667
$ws.Range("I8").Value = 170.66667
$ws.Range("K8").Value = 512.00001
$ws.Range("M8").Value = -373.00001

$ws.Range("H11").Value = 32.4
$ws.Range("I11").Value = 32.4
$ws.Range("K11").Value = 32.4
$ws.Range("M11").Value = 107.6

$ws.Range("H53").Value = 371.55554
$ws.Range("I53").Value = 389.85715
$ws.Range("J53").Value = 307.5
$ws.Range("K53").Value = 389.85715
$ws.Range("L53").Value = 307.5
$ws.Range("M53").Value = 247.14285
$ws.Range("N53").Value = -1581.5

$ws.Range("H116").Value = 2500
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = 1942
$ws.Range("N116").Value = -10384

$ws.Range("H138").Value = 6399.4
$ws.Range("J138").Value = 7749.25
$ws.Range("L138").Value = 23247.75
$ws.Range("N138").Value = -33527.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 30.066668
$ws.Range("I5").Value = 50.2
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 50.2
$ws.Range("L5").Value = 20
$ws.Range("M5").Value = 61.8
$ws.Range("N5").Value = -244

$ws.Range("H32").Value = 168.16667
$ws.Range("I32").Value = 168.16667
$ws.Range("K32").Value = 168.16667
$ws.Range("M32").Value = 118.83333

$ws.Range("H74").Value = 2277.75
$ws.Range("I74").Value = 2030.2727
$ws.Range("K74").Value = 2030.2727
$ws.Range("M74").Value = -1156.2727

$ws.Range("H77").Value = 2277.75
$ws.Range("I77").Value = 2030.2727
$ws.Range("K77").Value = 10151.3635
$ws.Range("M77").Value = -5783.363499999999

$ws.Range("H122").Value = 1141.2858
$ws.Range("I122").Value = 1039.8334
$ws.Range("K122").Value = 3119.5002
$ws.Range("M122").Value = -669.5001999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 30.066668
$ws.Range("I4").Value = 50.2
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 50.2
$ws.Range("L4").Value = 20
$ws.Range("M4").Value = 64.8
$ws.Range("N4").Value = -250

$ws.Range("H22").Value = 3142.3845
$ws.Range("I22").Value = 3008.4167
$ws.Range("K22").Value = 3008.4167
$ws.Range("M22").Value = -2835.4167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 12779266
$ws.Range("I6").Value = 13126674
$ws.Range("J6").Value = 10000000
$ws.Range("K6").Value = 13126674
$ws.Range("L6").Value = 10000000
$ws.Range("M6").Value = -13126561
$ws.Range("N6").Value = -10000226

$ws.Range("H10").Value = 1055.7059
$ws.Range("I10").Value = 620.1111
$ws.Range("K10").Value = 620.1111
$ws.Range("M10").Value = -481.1111

$ws.Range("H19").Value = 500.22726
$ws.Range("I19").Value = 253.41176
$ws.Range("J19").Value = 1339.4
$ws.Range("K19").Value = 253.41176
$ws.Range("L19").Value = 1339.4
$ws.Range("M19").Value = -83.41175999999999
$ws.Range("N19").Value = -1679.4

$ws.Range("H24").Value = 500.22726
$ws.Range("I24").Value = 253.41176
$ws.Range("J24").Value = 1339.4
$ws.Range("K24").Value = 253.41176
$ws.Range("L24").Value = 1339.4
$ws.Range("M24").Value = -83.41175999999999
$ws.Range("N24").Value = -1679.4

$ws.Range("H33").Value = 20479.834
$ws.Range("I33").Value = 720
$ws.Range("J33").Value = 59999.5
$ws.Range("K33").Value = 720
$ws.Range("L33").Value = 59999.5
$ws.Range("M33").Value = -341
$ws.Range("N33").Value = -60757.5

$ws.Range("H92").Value = 68999.75
$ws.Range("J92").Value = 68999.75
$ws.Range("L92").Value = 68999.75
$ws.Range("N92").Value = -73991.75

$ws.Range("H122").Value = 2162.5
$ws.Range("I122").Value = 1550
$ws.Range("K122").Value = 4650
$ws.Range("M122").Value = -2200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 637.5
$ws.Range("J12").Value = 1200
$ws.Range("L12").Value = 3600
$ws.Range("N12").Value = -3946

$ws.Range("H23").Value = 909.4286
$ws.Range("I23").Value = 89.666664
$ws.Range("K23").Value = 268.999992
$ws.Range("M23").Value = -33.99999200000002

$ws.Range("H25").Value = 37525
$ws.Range("I25").Value = 50
$ws.Range("J25").Value = 75000
$ws.Range("K25").Value = 150
$ws.Range("L25").Value = 225000
$ws.Range("M25").Value = 19
$ws.Range("N25").Value = -225338

$ws.Range("H30").Value = 37525
$ws.Range("I30").Value = 50
$ws.Range("J30").Value = 75000
$ws.Range("K30").Value = 150
$ws.Range("L30").Value = 225000
$ws.Range("M30").Value = -48
$ws.Range("N30").Value = -225204

$ws.Range("H80").Value = 2292.6667
$ws.Range("I80").Value = 937.5
$ws.Range("J80").Value = 5003
$ws.Range("K80").Value = 2812.5
$ws.Range("L80").Value = 15009
$ws.Range("M80").Value = -1876.5
$ws.Range("N80").Value = -16881

$ws.Range("H83").Value = 2292.6667
$ws.Range("I83").Value = 937.5
$ws.Range("J83").Value = 5003
$ws.Range("K83").Value = 8437.5
$ws.Range("L83").Value = 45027
$ws.Range("M83").Value = -3757.5
$ws.Range("N83").Value = -54387

$ws.Range("H97").Value = 866
$ws.Range("I97").Value = 866
$ws.Range("K97").Value = 2598
$ws.Range("M97").Value = -2102

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 51000
$ws.Range("J51").Value = 51000
$ws.Range("L51").Value = 51000
$ws.Range("N51").Value = -52018

$ws.Range("H97").Value = 1219.75
$ws.Range("I97").Value = 1214.5
$ws.Range("K97").Value = 1214.5
$ws.Range("M97").Value = -718.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H25").Value = 7000
$ws.Range("I25").Value = 7000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 7000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -6770
$ws.Range("N25").ClearContents()

$ws.Range("H40").Value = 45481.57
$ws.Range("I40").Value = 25673.2
$ws.Range("K40").Value = 25673.2
$ws.Range("M40").Value = -25537.2

$ws.Range("H43").Value = 349484.5
$ws.Range("J43").Value = 349484.5
$ws.Range("L43").Value = 349484.5
$ws.Range("N43").Value = -349870.5

$ws.Range("H46").Value = 798.5
$ws.Range("J46").Value = 798.2
$ws.Range("L46").Value = 798.2
$ws.Range("N46").Value = -1174.2

$ws.Range("H55").Value = 595
$ws.Range("I55").Value = 595
$ws.Range("K55").Value = 595
$ws.Range("M55").Value = -422

$ws.Range("H82").Value = 2000
$ws.Range("J82").Value = 2000
$ws.Range("L82").Value = 2000
$ws.Range("N82").Value = -2722

$ws.Range("H85").Value = 2000
$ws.Range("J85").Value = 2000
$ws.Range("L85").Value = 2000
$ws.Range("N85").Value = -4496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 192962
$ws.Range("I2").Value = 188888
$ws.Range("J2").Value = 194999
$ws.Range("K2").Value = 188888
$ws.Range("L2").Value = 194999
$ws.Range("M2").Value = -188776
$ws.Range("N2").Value = -195223

$ws.Range("H126").Value = 1775.8334
$ws.Range("I126").Value = 1530
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 4590
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -2120
$ws.Range("N126").Value = -13955

$ws.Range("H132").Value = 3601.25
$ws.Range("I132").Value = 3466.6667
$ws.Range("J132").Value = 4005
$ws.Range("K132").Value = 10400.0001
$ws.Range("L132").Value = 12015
$ws.Range("M132").Value = -7870.000100000001
$ws.Range("N132").Value = -17075
